# checkpoint antes de iniciar crossval
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Fill in the previously-empty "tree" accuracy columns (C and F) with the
# cross-validation checkpoint values.
$ws.Range("C4").Value = 0.375
$ws.Range("F4").Value = 0.33129999999999998

$ws.Range("C5").Value = 0.41249999999999998
$ws.Range("F5").Value = 0.28129999999999999

$ws.Range("C6").Value = 0.41880000000000001
$ws.Range("F6").Value = 0.28129999999999999

$ws.Range("C7").Value = 0.31059999999999999

$ws.Range("C8").Value = 0.25

$ws.Range("C9").Value = 0.43130000000000002

$ws.Range("C10").Value = 0.41880000000000001

$ws.Range("C11").Value = 0.39379999999999998

$ws.Range("C12").Value = 0.38129999999999997

$ws.Range("C13").Value = 0.32950000000000002

$ws.Range("C18").Value = 0.39129999999999998

$ws.Range("C19").Value = 0.4224

$ws.Range("C20").Value = 0.4037

$ws.Range("C21").Value = 0.32079999999999997

$ws.Range("C22").Value = 0.28129999999999999

$ws.Range("C23").Value = 0.3851

$ws.Range("C24").Value = 0.43669999999999998

$ws.Range("C25").Value = 0.45629999999999998

$ws.Range("C26").Value = 0.30249999999999999

$ws.Range("C27").Value = 0.38129999999999997

# Move the active selection to F7, matching the saved cursor position.
$ws.Activate()
$ws.Range("F7").Select()
